# Trade #14 (MarketMaking strategy, global trade #43) closed at
# 2026-02-18 00:10:58 as an early_exit with a small loss, and a brand new
# momentum trade (#72) was opened at 2026-02-18 00:10:52. This updates the
# rollup sheets (Summary / Strategy Status), the "All Trades" log, the
# per-strategy "MarketMaking" log (closing trade) and the per-strategy
# "momentum" log (new trade).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Summary sheet
# ---------------------------------------------------------------------
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B3").Value = 1499.67   # Current Capital
$summary.Range("B4").Value = 0.77      # Total P&L $
$summary.Range("B5").Value = 0.37      # Total P&L %
$summary.Range("B6").Value = 42        # Total Trades
$summary.Range("B8").Value = 16        # Losing Trades
$summary.Range("B9").Value = 54.76     # Win Rate %

# ---------------------------------------------------------------------
# Strategy Status sheet - MarketMaking row (row 6)
# ---------------------------------------------------------------------
$status = $wb.Worksheets.Item("Strategy Status")
$status.Range("C6").Value = 99.67    # Capital
$status.Range("D6").Value = 13       # Trades
$status.Range("E6").Value = -0.14    # P&L $
$status.Range("F6").Value = -0.33    # P&L %
$status.Range("G6").Value = 53.85    # Win Rate %

# ---------------------------------------------------------------------
# All Trades sheet - close out trade #43 (row 44, MarketMaking)
# ---------------------------------------------------------------------
$allTrades = $wb.Worksheets.Item("All Trades")
$allTrades.Cells.Item(44, 7).Value = 0.39          # G44 Exit Price
$allTrades.Cells.Item(44, 8).Value = "CLOSED"      # H44 Status
$allTrades.Cells.Item(44, 9).Value = -13.3333      # I44 P&L %
$allTrades.Cells.Item(44, 10).Value = -0.06        # J44 P&L $
$allTrades.Cells.Item(44, 11).Value = 99.67        # K44 Capital After
$allTrades.Cells.Item(44, 12).Value = "early_exit" # L44 Exit Reason
$allTrades.Cells.Item(44, 13).Value = 0.15         # M44 Duration (min)

# All Trades sheet - append new trade #72 (row 73, momentum, OPEN)
$allTrades.Cells.Item(73, 1).Value = 72
$allTrades.Cells.Item(73, 2).Value = "'2026-02-18"
$allTrades.Cells.Item(73, 3).Value = "'00:10:52"
$allTrades.Cells.Item(73, 4).Value = "momentum"
$allTrades.Cells.Item(73, 5).Value = "UP"
$allTrades.Cells.Item(73, 6).Value = 0.45
# G73 Exit Price left blank (trade is still open)
$allTrades.Cells.Item(73, 8).Value = "OPEN"
$allTrades.Cells.Item(73, 9).Value = 0
$allTrades.Cells.Item(73, 10).Value = 0
$allTrades.Cells.Item(73, 11).Value = 100
# L73 Exit Reason left blank (trade is still open)
$allTrades.Cells.Item(73, 13).Value = 0
$allTrades.Cells.Item(73, 14).Value = 0
$allTrades.Cells.Item(73, 15).Value = 0
$allTrades.Cells.Item(73, 16).Value = 0.9
$allTrades.Cells.Item(73, 17).Value = "Upward momentum: 21.687% over 10 samples"

# ---------------------------------------------------------------------
# momentum sheet - append new trade #72 (row 9, OPEN)
# ---------------------------------------------------------------------
$momentum = $wb.Worksheets.Item("momentum")
$momentum.Cells.Item(9, 1).Value = 72
$momentum.Cells.Item(9, 2).Value = "'2026-02-18"
$momentum.Cells.Item(9, 3).Value = "'00:10:52"
$momentum.Cells.Item(9, 4).Value = "momentum"
$momentum.Cells.Item(9, 5).Value = "UP"
$momentum.Cells.Item(9, 6).Value = 0.45
# G9 Exit Price left blank (trade is still open)
$momentum.Cells.Item(9, 8).Value = "OPEN"
$momentum.Cells.Item(9, 9).Value = 0
$momentum.Cells.Item(9, 10).Value = 0
$momentum.Cells.Item(9, 11).Value = 100
$momentum.Cells.Item(9, 12).Value = 0
$momentum.Cells.Item(9, 13).Value = 0
$momentum.Cells.Item(9, 14).Value = 0.9
$momentum.Cells.Item(9, 15).Value = "Upward momentum: 21.687% over 10 samples"
# P9 Exit Reason left blank (trade is still open)
$momentum.Cells.Item(9, 17).Value = 0

# ---------------------------------------------------------------------
# MarketMaking sheet - close out trade #43 (row 15)
# ---------------------------------------------------------------------
$marketMaking = $wb.Worksheets.Item("MarketMaking")
$marketMaking.Cells.Item(15, 7).Value = 0.39          # G15 Exit Price
$marketMaking.Cells.Item(15, 8).Value = "CLOSED"      # H15 Status
$marketMaking.Cells.Item(15, 9).Value = -13.3333      # I15 P&L %
$marketMaking.Cells.Item(15, 10).Value = -0.06        # J15 P&L $
$marketMaking.Cells.Item(15, 11).Value = 99.67        # K15 Capital After
$marketMaking.Cells.Item(15, 16).Value = "early_exit" # P15 Exit Reason
$marketMaking.Cells.Item(15, 17).Value = 0.15         # Q15 Duration (min)
